$wb = $excel.ActiveWorkbook

# This script applies the numeric corrections captured in the commit diff
# for the "Bahamut_Profits" workbook snapshot. Each worksheet corresponds to
# a crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) and stores leve
# profitability data in a Table (A1:N141) with columns H..N holding the
# market-price/profit calculations that were refreshed by the scheduled runner.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 3231.1
$ws.Range("I96").Value = 3330.2144
$ws.Range("J96").Value = 2999.8333
$ws.Range("K96").Value = 9990.643199999999
$ws.Range("L96").Value = 8999.499899999999
$ws.Range("M96").Value = -8617.643199999999
$ws.Range("N96").Value = -11745.4999
$ws.Range("H101").Value = 125583
$ws.Range("I101").Value = 447.33334
$ws.Range("K101").Value = 1342.00002
$ws.Range("M101").Value = 279.9999800000001
$ws.Range("H107").Value = 2500332.8
$ws.Range("I107").Value = 5000090.5
$ws.Range("J107").Value = 575
$ws.Range("K107").Value = 5000090.5
$ws.Range("L107").Value = 575
$ws.Range("M107").Value = -4998170.5
$ws.Range("N107").Value = -4415
$ws.Range("H137").Value = 1089.1666
$ws.Range("I137").Value = 911
$ws.Range("J137").Value = 1980
$ws.Range("K137").Value = 2733
$ws.Range("L137").Value = 5940
$ws.Range("M137").Value = -183
$ws.Range("N137").Value = -11040

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6974.4414
$ws.Range("I32").Value = 6156.185
$ws.Range("J32").Value = 10130.571
$ws.Range("K32").Value = 6156.185
$ws.Range("L32").Value = 10130.571
$ws.Range("M32").Value = -5869.185
$ws.Range("N32").Value = -10704.571
$ws.Range("H61").Value = 2960
$ws.Range("I61").Value = 2933.3333
$ws.Range("K61").Value = 2933.3333
$ws.Range("M61").Value = -2721.3333
$ws.Range("H74").Value = 887.5925999999999
$ws.Range("I74").Value = 803.3
$ws.Range("J74").Value = 1128.4286
$ws.Range("K74").Value = 803.3
$ws.Range("L74").Value = 1128.4286
$ws.Range("M74").Value = 70.70000000000005
$ws.Range("N74").Value = -2876.4286
$ws.Range("H77").Value = 887.5925999999999
$ws.Range("I77").Value = 803.3
$ws.Range("J77").Value = 1128.4286
$ws.Range("K77").Value = 4016.5
$ws.Range("L77").Value = 5642.143
$ws.Range("M77").Value = 351.5
$ws.Range("N77").Value = -14378.143
$ws.Range("H92").Value = 27750
$ws.Range("J92").Value = 27750
$ws.Range("L92").Value = 27750
$ws.Range("N92").Value = -32742
$ws.Range("H102").Value = 4265.3335
$ws.Range("I102").Value = 4246.25
$ws.Range("J102").Value = 4287.143
$ws.Range("K102").Value = 4246.25
$ws.Range("L102").Value = 4287.143
$ws.Range("M102").Value = -2624.25
$ws.Range("N102").Value = -7531.143
$ws.Range("H122").Value = 1355.5714
$ws.Range("I122").Value = 1355.5714
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4066.7142
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1616.7142
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 2960
$ws.Range("I136").Value = 2933.3333
$ws.Range("K136").Value = 8799.999899999999
$ws.Range("M136").Value = -6249.999899999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2918.75
$ws.Range("I86").Value = 2192.8572
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 2192.8572
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -1069.8572
$ws.Range("N86").Value = -10246
$ws.Range("H89").Value = 2918.75
$ws.Range("I89").Value = 2192.8572
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 10964.286
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -5348.286
$ws.Range("N89").Value = -51232
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H94").Value = 851.2727
$ws.Range("I94").Value = 499.8
$ws.Range("J94").Value = 1144.1666
$ws.Range("K94").Value = 499.8
$ws.Range("L94").Value = 1144.1666
$ws.Range("M94").Value = -48.80000000000001
$ws.Range("N94").Value = -2046.1666

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 19865
$ws.Range("J70").Value = 19865
$ws.Range("L70").Value = 19865
$ws.Range("N70").Value = -20495
$ws.Range("H73").Value = 19865
$ws.Range("J73").Value = 19865
$ws.Range("L73").Value = 19865
$ws.Range("N73").Value = -22049
$ws.Range("H80").Value = 20300
$ws.Range("J80").Value = 20300
$ws.Range("L80").Value = 20300
$ws.Range("N80").Value = -22546
$ws.Range("H83").Value = 20300
$ws.Range("J83").Value = 20300
$ws.Range("L83").Value = 60900
$ws.Range("N83").Value = -72132
$ws.Range("H92").Value = 21316.834
$ws.Range("J92").Value = 21316.834
$ws.Range("L92").Value = 21316.834
$ws.Range("N92").Value = -26308.834
$ws.Range("H97").Value = 12598.5
$ws.Range("J97").Value = 12598.5
$ws.Range("L97").Value = 12598.5
$ws.Range("N97").Value = -14580.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 5347811.5
$ws.Range("J2").Value = 332.2857
$ws.Range("L2").Value = 1993.7142
$ws.Range("N2").Value = -2219.7142
$ws.Range("H7").Value = 355
$ws.Range("I7").Value = 250
$ws.Range("K7").Value = 750
$ws.Range("M7").Value = -638

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10576
$ws.Range("H70").Value = 4395.4194
$ws.Range("J70").Value = 4828.385
$ws.Range("L70").Value = 4828.385
$ws.Range("N70").Value = -5368.385
$ws.Range("H73").Value = 4395.4194
$ws.Range("J73").Value = 4828.385
$ws.Range("L73").Value = 4828.385
$ws.Range("N73").Value = -6700.385
$ws.Range("H81").Value = 10000
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -11996
$ws.Range("H84").Value = 10000
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -39984
$ws.Range("H92").Value = 6166.6665
$ws.Range("J92").Value = 6166.6665
$ws.Range("L92").Value = 6166.6665
$ws.Range("N92").Value = -9910.666499999999
$ws.Range("H132").Value = 3905.1428
$ws.Range("I132").Value = 3106
$ws.Range("J132").Value = 4970.6665
$ws.Range("K132").Value = 9318
$ws.Range("L132").Value = 14911.9995
$ws.Range("M132").Value = -6788
$ws.Range("N132").Value = -19971.9995

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 322.66666
$ws.Range("I46").Value = 367
$ws.Range("J46").Value = 234
$ws.Range("K46").Value = 367
$ws.Range("L46").Value = 234
$ws.Range("M46").Value = -179
$ws.Range("N46").Value = -610
$ws.Range("H92").Value = 19425.945
$ws.Range("J92").Value = 19425.945
$ws.Range("L92").Value = 19425.945
$ws.Range("N92").Value = -24417.945
$ws.Range("H94").Value = 40000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 40000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 40000
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -41352
$ws.Range("H100").Value = 18520252
$ws.Range("I100").Value = 22223822
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 22223822
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -22223281
$ws.Range("N100").Value = -3482
$ws.Range("H136").Value = 5210.8
$ws.Range("I136").Value = 1301.1428
$ws.Range("J136").Value = 14333.333
$ws.Range("K136").Value = 3903.4284
$ws.Range("L136").Value = 42999.999
$ws.Range("M136").Value = -1353.4284
$ws.Range("N136").Value = -48099.999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 442.76923
$ws.Range("I113").Value = 449.55554
$ws.Range("J113").Value = 427.5
$ws.Range("K113").Value = 1348.66662
$ws.Range("L113").Value = 1282.5
$ws.Range("M113").Value = 821.33338
$ws.Range("N113").Value = -5622.5
$ws.Range("H136").Value = 3914.2
$ws.Range("I136").Value = 4271.933
$ws.Range("K136").Value = 12815.799
$ws.Range("M136").Value = -10265.799
